$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item("TextBox 2")
$tf = $shp.TextFrame
$tr = $tf.TextRange

$quoteL = [char]8220
$quoteR = [char]8221

$fullText = $tr.Text

# --- Paragraph: "Compute cost of all the neighbors ..." -----------------
# Collapse the many runs ("is ", "14, ", "9, ", "and ", "7, ", "respectively.")
# back into a single run with unified text (no wording/number change).
$computeNewText = "Compute cost of all the neighbors of the starting note (here: A). For instance, the cost of reaching B, C, and D from node A is 14, 9, and 7, respectively."
$computeMarker = "Compute cost of all the neighbors"
$computeEndMarker = "respectively."
$computeStart = $fullText.IndexOf($computeMarker) + 1
$computeEnd0 = $fullText.IndexOf($computeEndMarker) + $computeEndMarker.Length
$computeLen = $computeEnd0 - ($computeStart - 1)

$computeRange = $tr.Characters($computeStart, $computeLen)
$computeRange.Text = $computeNewText

# Recompute full text (length changed after the collapse above).
$fullText = $tr.Text

# --- Paragraph: "Assign cost "infinity" ..." -----------------------------
# Split the single run into five runs, changing the letter D -> E along the
# way: "(here: D and F)." -> "(here: E and F)."
$assignMarker = "Assign cost"
$assignStart = $fullText.IndexOf($assignMarker) + 1

$seg1 = "Assign cost " + $quoteL + "infinity" + $quoteR + " to all remaining nodes in the graph -- that is, all nodes that are not direct neighbors of the starting node A (here"
$seg2 = ": "
$seg3 = "E"
$seg4 = " "
$seg5 = "and F)."

$pos = $assignStart

$r1 = $tr.Characters($pos, $seg1.Length)
$r1.Text = $seg1
$pos += $seg1.Length

$r2 = $tr.Characters($pos, $seg2.Length)
$r2.Text = $seg2
$pos += $seg2.Length

$r3 = $tr.Characters($pos, $seg3.Length)
$r3.Text = $seg3
$pos += $seg3.Length

$r4 = $tr.Characters($pos, $seg4.Length)
$r4.Text = $seg4
$pos += $seg4.Length

$r5 = $tr.Characters($pos, $seg5.Length)
$r5.Text = $seg5
